$wb = $excel.ActiveWorkbook

# Sheet "OFF": update row 2 (Home) values
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 176
$wsOff.Range("C2").Value = 118
$wsOff.Range("D2").Value = 41
$wsOff.Range("E2").Value = 23

# Sheet "DEF": update row 2 (Home) values
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 219
$wsDef.Range("C2").Value = 149
$wsDef.Range("D2").Value = 56
$wsDef.Range("E2").Value = 30
